function Set-TextValue {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 'D2' '62.731.67'
Set-TextValue $ws 'E2' '  -0.74%  '
Set-TextValue $ws 'D3' '3.020.54'
Set-TextValue $ws 'E3' '  -0.91%  '
Set-TextValue $ws 'E4' '  -0.05%  '
Set-TextValue $ws 'D5' '587.13'
Set-TextValue $ws 'E5' '  -0.62%  '
Set-TextValue $ws 'D6' '148.57'
Set-TextValue $ws 'E6' '  -2.54%  '
Set-TextValue $ws 'E7' '  -0.05%  '
Set-TextValue $ws 'D8' '0.529'
Set-TextValue $ws 'E8' '  -2.27%  '
Set-TextValue $ws 'D9' '3.019.90'
Set-TextValue $ws 'E9' '  -0.89%  '
Set-TextValue $ws 'D10' '0.150'
Set-TextValue $ws 'E10' '  -2.47%  '
Set-TextValue $ws 'D11' '5.85'
Set-TextValue $ws 'E11' '  +0.70%  '
Set-TextValue $ws 'E12' '  +2.84%  '
Set-TextValue $ws 'D13' '0.0000232'
Set-TextValue $ws 'E13' '  -1.60%  '
Set-TextValue $ws 'D14' '34.93'
Set-TextValue $ws 'E14' '  -4.41%  '
Set-TextValue $ws 'E15' '  +2.14%  '
Set-TextValue $ws 'D16' '3.520.16'
Set-TextValue $ws 'E16' '  -0.92%  '
Set-TextValue $ws 'D17' '7.13'
Set-TextValue $ws 'E17' '  -0.45%  '
Set-TextValue $ws 'D18' '62.687.26'
Set-TextValue $ws 'E18' '  -0.75%  '
Set-TextValue $ws 'D19' '3.017.17'
Set-TextValue $ws 'E19' '  -1.06%  '
Set-TextValue $ws 'D20' '461.04'
Set-TextValue $ws 'E20' '  -4.10%  '
Set-TextValue $ws 'D21' '14.05'
Set-TextValue $ws 'E21' '  -1.53%  '
Set-TextValue $ws 'D22' '0.693'
Set-TextValue $ws 'E22' '  -1.77%  '
Set-TextValue $ws 'D23' '7.47'
Set-TextValue $ws 'E23' '  -0.73%  '
Set-TextValue $ws 'D24' '81.85'
Set-TextValue $ws 'E24' '  -0.29%  '
Set-TextValue $ws 'D25' '2.23'
Set-TextValue $ws 'E25' '  -8.22%  '
Set-TextValue $ws 'D26' '12.38'
Set-TextValue $ws 'E26' '  -2.99%  '
Set-TextValue $ws 'D27' '10.06'
Set-TextValue $ws 'E27' '  -6.04%  '
Set-TextValue $ws 'E28' '  +0.27%  '
Set-TextValue $ws 'B29' 'PancakeSwap'
Set-TextValue $ws 'C29' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws 'D29' '2.66'
Set-TextValue $ws 'E29' '  -0.75%  '
Set-TextValue $ws 'B30' 'FirstDigitalUSD'
Set-TextValue $ws 'C30' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws 'D30' '0.999'
Set-TextValue $ws 'E30' '  -0.09%  '
Set-TextValue $ws 'D31' '7.07'
Set-TextValue $ws 'E31' '  -4.40%  '
Set-TextValue $ws 'D32' '2.12'
Set-TextValue $ws 'E32' '  -4.20%  '
Set-TextValue $ws 'D33' '28.17'
Set-TextValue $ws 'E33' '  +2.51%  '
Set-TextValue $ws 'D34' '0.110'
Set-TextValue $ws 'E34' '  -1.43%  '
Set-TextValue $ws 'D35' '0.0₃0824'
Set-TextValue $ws 'E35' '  +0.40%  '
Set-TextValue $ws 'E36' '  -2.60%  '
Set-TextValue $ws 'D37' '5.81'
Set-TextValue $ws 'E37' '  -2.17%  '
Set-TextValue $ws 'D38' '2.14'
Set-TextValue $ws 'E38' '  -3.97%  '
Set-TextValue $ws 'D39' '50.46'
Set-TextValue $ws 'E39' '  +0.16%  '
Set-TextValue $ws 'D40' '9.18'
Set-TextValue $ws 'E40' '  -1.02%  '
Set-TextValue $ws 'D41' '2.94'
Set-TextValue $ws 'E41' '  -10.67%  '
Set-TextValue $ws 'D42' '0.124'
Set-TextValue $ws 'E42' '  +9.04%  '
Set-TextValue $ws 'D43' '396.26'
Set-TextValue $ws 'E43' '  -9.20%  '
Set-TextValue $ws 'D44' '0.0361'
Set-TextValue $ws 'E44' '  -0.63%  '
Set-TextValue $ws 'D45' '0.271'
Set-TextValue $ws 'E45' '  -6.33%  '
Set-TextValue $ws 'D46' '2.741.27'
Set-TextValue $ws 'E46' '  -3.10%  '
Set-TextValue $ws 'D47' '37.40'
Set-TextValue $ws 'E47' '  -2.21%  '
Set-TextValue $ws 'D48' '129.35'
Set-TextValue $ws 'E48' '  -0.61%  '
Set-TextValue $ws 'E50' '  +0.08%  '
Set-TextValue $ws 'D51' '2.21'
Set-TextValue $ws 'E51' '  -0.56%  '
